$d = $word.ActiveDocument

# 1) Trim the "Arbetsmarknadsutbildning" sentence: remove "var mindre givande och "
#    so "... tyvärr var mindre givande och avbröts ..." becomes "... tyvärr avbröts ...".
$d.Content.Find.Execute(
    "var mindre givande och avbröts",  # FindText
    $true,                             # MatchCase
    $false,                            # MatchWholeWord
    $false,                            # MatchWildcards
    $false,                            # MatchSoundsLike
    $false,                            # MatchAllWordForms
    $true,                             # Forward
    1,                                 # Wrap (wdFindContinue)
    $false,                            # Format
    "avbröts",                         # ReplaceWith
    2                                  # Replace (wdReplaceAll)
) | Out-Null

# 2) Center-align the final paragraph (name / email / phone contact line).
$last = $d.Paragraphs.Last
$last.Alignment = 1
